$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 777.73914
$ws.Cells.Item(15, 9).Value = 777.73914
$ws.Cells.Item(15, 11).Value = 2333.21742
$ws.Cells.Item(15, 13).Value = -2164.21742
$ws.Cells.Item(32, 8).Value = 2159.3333
$ws.Cells.Item(32, 10).Value = 2299.182
$ws.Cells.Item(32, 12).Value = 2299.182
$ws.Cells.Item(32, 14).Value = -2951.182
$ws.Cells.Item(81, 8).Value = 750018750
$ws.Cells.Item(81, 10).Value = 750018750
$ws.Cells.Item(81, 12).Value = 750018750
$ws.Cells.Item(81, 14).Value = -750020746
$ws.Cells.Item(84, 8).Value = 750018750
$ws.Cells.Item(84, 10).Value = 750018750
$ws.Cells.Item(84, 12).Value = 2250056250
$ws.Cells.Item(84, 14).Value = -2250066234
$ws.Cells.Item(94, 8).Value = 1356.8572
$ws.Cells.Item(94, 9).Value = 1356.8572
$ws.Cells.Item(94, 11).Value = 1356.8572
$ws.Cells.Item(94, 13).Value = -905.8571999999999
$ws.Cells.Item(129, 8).Value = 1723.4615
$ws.Cells.Item(129, 9).Value = 1162.6666
$ws.Cells.Item(129, 11).Value = 3487.9998
$ws.Cells.Item(129, 13).Value = 1512.0002
$ws.Cells.Item(138, 8).Value = 2551.9246
$ws.Cells.Item(138, 10).Value = 3156.5
$ws.Cells.Item(138, 12).Value = 9469.5
$ws.Cells.Item(138, 14).Value = -19749.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2849.4583
$ws.Cells.Item(32, 9).Value = 2659.6191
$ws.Cells.Item(32, 11).Value = 2659.6191
$ws.Cells.Item(32, 13).Value = -2372.6191
$ws.Cells.Item(61, 8).Value = 2397.9473
$ws.Cells.Item(61, 9).Value = 2430.1428
$ws.Cells.Item(61, 11).Value = 2430.1428
$ws.Cells.Item(61, 13).Value = -2218.1428
$ws.Cells.Item(74, 8).Value = 975.41174
$ws.Cells.Item(74, 9).Value = 714.5161000000001
$ws.Cells.Item(74, 10).Value = 3671.3333
$ws.Cells.Item(74, 11).Value = 714.5161000000001
$ws.Cells.Item(74, 12).Value = 3671.3333
$ws.Cells.Item(74, 13).Value = 159.4838999999999
$ws.Cells.Item(74, 14).Value = -5419.3333
$ws.Cells.Item(77, 8).Value = 975.41174
$ws.Cells.Item(77, 9).Value = 714.5161000000001
$ws.Cells.Item(77, 10).Value = 3671.3333
$ws.Cells.Item(77, 11).Value = 3572.5805
$ws.Cells.Item(77, 12).Value = 18356.6665
$ws.Cells.Item(77, 13).Value = 795.4195
$ws.Cells.Item(77, 14).Value = -27092.6665
$ws.Cells.Item(132, 8).Value = 2155.4375
$ws.Cells.Item(132, 9).Value = 2132.4666
$ws.Cells.Item(132, 10).Value = 2500
$ws.Cells.Item(132, 11).Value = 6397.399800000001
$ws.Cells.Item(132, 12).Value = 7500
$ws.Cells.Item(132, 13).Value = -3867.399800000001
$ws.Cells.Item(132, 14).Value = -12560
$ws.Cells.Item(136, 8).Value = 2397.9473
$ws.Cells.Item(136, 9).Value = 2430.1428
$ws.Cells.Item(136, 11).Value = 7290.428400000001
$ws.Cells.Item(136, 13).Value = -4740.428400000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 1588.1666
$ws.Cells.Item(22, 9).Value = 1505.8
$ws.Cells.Item(22, 11).Value = 1505.8
$ws.Cells.Item(22, 13).Value = -1332.8
$ws.Cells.Item(99, 8).Value = 2768.4375
$ws.Cells.Item(99, 9).Value = 1183.3334
$ws.Cells.Item(99, 10).Value = 3719.5
$ws.Cells.Item(99, 11).Value = 1183.3334
$ws.Cells.Item(99, 12).Value = 3719.5
$ws.Cells.Item(99, 13).Value = 314.6666
$ws.Cells.Item(99, 14).Value = -6715.5
$ws.Cells.Item(134, 8).Value = 1706.0625
$ws.Cells.Item(134, 9).Value = 1461.3846
$ws.Cells.Item(134, 11).Value = 4384.1538
$ws.Cells.Item(134, 13).Value = -1849.1538
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1614.5385
$ws.Cells.Item(58, 9).Value = 1713.5
$ws.Cells.Item(58, 10).Value = 1499.0834
$ws.Cells.Item(58, 11).Value = 1713.5
$ws.Cells.Item(58, 12).Value = 1499.0834
$ws.Cells.Item(58, 13).Value = -1510.5
$ws.Cells.Item(58, 14).Value = -1905.0834
$ws.Cells.Item(99, 8).Value = 2978.2856
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 14).Value = ""
$ws.Cells.Item(103, 8).Value = 13666.667
$ws.Cells.Item(103, 9).Value = 13666.667
$ws.Cells.Item(103, 11).Value = 13666.667
$ws.Cells.Item(103, 13).Value = -12494.667
$ws.Cells.Item(126, 8).Value = 2978.2856
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 14).Value = ""
$ws.Cells.Item(132, 8).Value = 2248.2354
$ws.Cells.Item(132, 10).Value = 1499.6666
$ws.Cells.Item(132, 12).Value = 4498.9998
$ws.Cells.Item(132, 14).Value = -9558.9998
$ws.Cells.Item(134, 8).Value = 2579.2273
$ws.Cells.Item(134, 9).Value = 2612.15
$ws.Cells.Item(134, 11).Value = 7836.450000000001
$ws.Cells.Item(134, 13).Value = -5301.450000000001
$ws.Cells.Item(136, 8).Value = 1614.5385
$ws.Cells.Item(136, 9).Value = 1713.5
$ws.Cells.Item(136, 10).Value = 1499.0834
$ws.Cells.Item(136, 11).Value = 5140.5
$ws.Cells.Item(136, 12).Value = 4497.2502
$ws.Cells.Item(136, 13).Value = -2590.5
$ws.Cells.Item(136, 14).Value = -9597.2502
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 105.15385
$ws.Cells.Item(33, 9).Value = 52.166668
$ws.Cells.Item(33, 10).Value = 150.57143
$ws.Cells.Item(33, 11).Value = 313.000008
$ws.Cells.Item(33, 12).Value = 903.42858
$ws.Cells.Item(33, 13).Value = -30.00000799999998
$ws.Cells.Item(33, 14).Value = -1469.42858
$ws.Cells.Item(121, 8).Value = 857.93335
$ws.Cells.Item(121, 9).Value = 422.16666
$ws.Cells.Item(121, 10).Value = 1148.4445
$ws.Cells.Item(121, 11).Value = 1266.49998
$ws.Cells.Item(121, 12).Value = 3445.3335
$ws.Cells.Item(121, 13).Value = 43.50001999999995
$ws.Cells.Item(121, 14).Value = -6065.333500000001
$ws.Cells.Item(129, 8).Value = 2021.6
$ws.Cells.Item(129, 10).Value = 2434.875
$ws.Cells.Item(129, 12).Value = 7304.625
$ws.Cells.Item(129, 14).Value = -17304.625
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(59, 8).Value = 49750
$ws.Cells.Item(59, 10).Value = 49500
$ws.Cells.Item(59, 12).Value = 49500
$ws.Cells.Item(59, 14).Value = -50666
$ws.Cells.Item(80, 8).Value = 5750.7856
$ws.Cells.Item(80, 9).Value = 4886.1665
$ws.Cells.Item(80, 10).Value = 6399.25
$ws.Cells.Item(80, 11).Value = 4886.1665
$ws.Cells.Item(80, 12).Value = 6399.25
$ws.Cells.Item(80, 13).Value = -3888.1665
$ws.Cells.Item(80, 14).Value = -8395.25
$ws.Cells.Item(83, 8).Value = 5750.7856
$ws.Cells.Item(83, 9).Value = 4886.1665
$ws.Cells.Item(83, 10).Value = 6399.25
$ws.Cells.Item(83, 11).Value = 24430.8325
$ws.Cells.Item(83, 12).Value = 31996.25
$ws.Cells.Item(83, 13).Value = -19438.8325
$ws.Cells.Item(83, 14).Value = -41980.25
$ws.Cells.Item(122, 8).Value = 5832
$ws.Cells.Item(122, 9).Value = 5330.6665
$ws.Cells.Item(122, 10).Value = 6333.3335
$ws.Cells.Item(122, 11).Value = 15991.9995
$ws.Cells.Item(122, 12).Value = 19000.0005
$ws.Cells.Item(122, 13).Value = -13541.9995
$ws.Cells.Item(122, 14).Value = -23900.0005
$ws.Cells.Item(132, 8).Value = 2848.2354
$ws.Cells.Item(132, 9).Value = 1676.7273
$ws.Cells.Item(132, 10).Value = 4996
$ws.Cells.Item(132, 11).Value = 5030.1819
$ws.Cells.Item(132, 12).Value = 14988
$ws.Cells.Item(132, 13).Value = -2500.1819
$ws.Cells.Item(132, 14).Value = -20048
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2513
$ws.Cells.Item(22, 9).Value = 2641.4167
$ws.Cells.Item(22, 10).Value = 1999.3334
$ws.Cells.Item(22, 11).Value = 2641.4167
$ws.Cells.Item(22, 12).Value = 1999.3334
$ws.Cells.Item(22, 13).Value = -2346.4167
$ws.Cells.Item(22, 14).Value = -2589.3334
$ws.Cells.Item(27, 8).Value = 2513
$ws.Cells.Item(27, 9).Value = 2641.4167
$ws.Cells.Item(27, 10).Value = 1999.3334
$ws.Cells.Item(27, 11).Value = 2641.4167
$ws.Cells.Item(27, 12).Value = 1999.3334
$ws.Cells.Item(27, 13).Value = -2534.4167
$ws.Cells.Item(27, 14).Value = -2213.3334
$ws.Cells.Item(46, 8).Value = 1624.6364
$ws.Cells.Item(46, 10).Value = 1534.375
$ws.Cells.Item(46, 12).Value = 1534.375
$ws.Cells.Item(46, 14).Value = -1910.375
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 13645.5
$ws.Cells.Item(62, 10).Value = 14023.2
$ws.Cells.Item(62, 12).Value = 14023.2
$ws.Cells.Item(62, 14).Value = -15271.2
$ws.Cells.Item(65, 8).Value = 13645.5
$ws.Cells.Item(65, 10).Value = 14023.2
$ws.Cells.Item(65, 12).Value = 70116
$ws.Cells.Item(65, 14).Value = -76356
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 14).Value = ""
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).Value = ""
$ws.Cells.Item(96, 8).Value = 4501.5
$ws.Cells.Item(96, 9).Value = 4501.5
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 11).Value = 4501.5
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 13).Value = -3128.5
$ws.Cells.Item(96, 14).Value = ""
$ws.Cells.Item(107, 8).Value = 1635.8182
$ws.Cells.Item(107, 10).Value = 1845
$ws.Cells.Item(107, 12).Value = 5535
$ws.Cells.Item(107, 14).Value = -9375
